$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1880566.6
$ws.Range("I19").Value = 4386594
$ws.Range("K19").Value = 4386594
$ws.Range("M19").Value = -4386419
# Row 28
$ws.Range("H28").Value = 752.64703
$ws.Range("I28").Value = 742.5
$ws.Range("J28").Value = 800
$ws.Range("K28").Value = 742.5
$ws.Range("L28").Value = 800
$ws.Range("M28").Value = -257.5
$ws.Range("N28").Value = -1770
# Row 40
$ws.Range("H40").Value = 1945
$ws.Range("J40").Value = 2193.3333
$ws.Range("L40").Value = 2193.3333
$ws.Range("N40").Value = -2543.3333
# Row 52
$ws.Range("H52").Value = 101700
$ws.Range("J52").Value = 101700
$ws.Range("L52").Value = 305100
$ws.Range("N52").Value = -305420
# Row 80
$ws.Range("H80").Value = 747.381
$ws.Range("J80").Value = 1322.5
$ws.Range("L80").Value = 3967.5
$ws.Range("N80").Value = -5963.5
# Row 83
$ws.Range("H83").Value = 747.381
$ws.Range("J83").Value = 1322.5
$ws.Range("L83").Value = 11902.5
$ws.Range("N83").Value = -21886.5
# Row 112
$ws.Range("H112").Value = 1260.0702
$ws.Range("J112").Value = 1328.7548
$ws.Range("L112").Value = 3986.2644
$ws.Range("N112").Value = -6202.2644
# Row 116
$ws.Range("H116").Value = 375766.4
$ws.Range("I116").Value = 835722.9399999999
$ws.Range("J116").Value = 7801.2
$ws.Range("K116").Value = 835722.9399999999
$ws.Range("L116").Value = 7801.2
$ws.Range("M116").Value = -832280.9399999999
$ws.Range("N116").Value = -14685.2
# Row 132
$ws.Range("H132").Value = 34484880
$ws.Range("I132").Value = 37038690
$ws.Range("K132").Value = 111116070
$ws.Range("M132").Value = -111113540

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 571
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 713
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 713
$ws.Range("M2").Value = -387
$ws.Range("N2").Value = -939
# Row 5
$ws.Range("H5").Value = 261.75
$ws.Range("I5").Value = 283
$ws.Range("K5").Value = 283
$ws.Range("M5").Value = -171
# Row 32
$ws.Range("H32").Value = 6720.4736
$ws.Range("I32").Value = 5481.1514
$ws.Range("K32").Value = 5481.1514
$ws.Range("M32").Value = -5194.1514
# Row 45
$ws.Range("H45").Value = 1278.8334
$ws.Range("I45").Value = 1278.8334
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1278.8334
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -901.8334
$ws.Range("N45").ClearContents()
# Row 80
$ws.Range("H80").Value = 37818.668
$ws.Range("J80").Value = 37818.668
$ws.Range("L80").Value = 37818.668
$ws.Range("N80").Value = -39814.668
# Row 83
$ws.Range("H83").Value = 37818.668
$ws.Range("J83").Value = 37818.668
$ws.Range("L83").Value = 113456.004
$ws.Range("N83").Value = -123440.004
# Row 109
$ws.Range("H109").Value = 30900
$ws.Range("J109").Value = 30900
$ws.Range("L109").Value = 30900
$ws.Range("N109").Value = -33674
# Row 116
$ws.Range("H116").Value = 571
$ws.Range("I116").Value = 500
$ws.Range("J116").Value = 713
$ws.Range("K116").Value = 500
$ws.Range("L116").Value = 713
$ws.Range("M116").Value = 1794
$ws.Range("N116").Value = -5301
# Row 132
$ws.Range("H132").Value = 6027.75
$ws.Range("J132").Value = 7466.3335
$ws.Range("L132").Value = 22399.0005
$ws.Range("N132").Value = -27459.0005
# Row 134
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
# Row 135
$ws.Range("H135").Value = 49929
$ws.Range("J135").Value = 49929
$ws.Range("L135").Value = 49929
$ws.Range("N135").Value = -60069
# Row 137
$ws.Range("H137").Value = 38793.6
$ws.Range("J137").Value = 41042
$ws.Range("L137").Value = 41042
$ws.Range("N137").Value = -51242

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 571
$ws.Range("I3").Value = 500
$ws.Range("J3").Value = 713
$ws.Range("K3").Value = 500
$ws.Range("L3").Value = 713
$ws.Range("M3").Value = -386
$ws.Range("N3").Value = -941
# Row 4
$ws.Range("H4").Value = 261.75
$ws.Range("I4").Value = 283
$ws.Range("K4").Value = 283
$ws.Range("M4").Value = -168
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
# Row 38
$ws.Range("H38").Value = 14269.75
$ws.Range("J38").Value = 14269.75
$ws.Range("L38").Value = 14269.75
$ws.Range("N38").Value = -15101.75
# Row 99
$ws.Range("H99").Value = 5344.5557
$ws.Range("I99").Value = 1490
$ws.Range("J99").Value = 5826.375
$ws.Range("K99").Value = 1490
$ws.Range("L99").Value = 5826.375
$ws.Range("M99").Value = 8
$ws.Range("N99").Value = -8822.375
# Row 107
$ws.Range("H107").Value = 2441.0715
$ws.Range("I107").Value = 1785
$ws.Range("K107").Value = 1785
$ws.Range("M107").Value = 135
# Row 114
$ws.Range("H114").Value = 30342
$ws.Range("J114").Value = 30342
$ws.Range("L114").Value = 30342
$ws.Range("N114").Value = -39020
# Row 134
$ws.Range("H134").Value = 3284.238
$ws.Range("I134").Value = 2516.8125
$ws.Range("J134").Value = 5740
$ws.Range("K134").Value = 7550.4375
$ws.Range("L134").Value = 17220
$ws.Range("M134").Value = -5015.4375
$ws.Range("N134").Value = -22290

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 358.94446
$ws.Range("I7").Value = 360.1111
$ws.Range("J7").Value = 357.77777
$ws.Range("K7").Value = 360.1111
$ws.Range("L7").Value = 357.77777
$ws.Range("M7").Value = -247.1111
$ws.Range("N7").Value = -583.7777699999999
# Row 58
$ws.Range("H58").Value = 3113.492
$ws.Range("I58").Value = 1783.02
$ws.Range("K58").Value = 1783.02
$ws.Range("M58").Value = -1580.02
# Row 105
$ws.Range("H105").Value = 1358.9706
$ws.Range("I105").Value = 1017.3333
$ws.Range("J105").Value = 2178.9
$ws.Range("K105").Value = 1017.3333
$ws.Range("L105").Value = 2178.9
$ws.Range("M105").Value = 729.6667
$ws.Range("N105").Value = -5672.9
# Row 107
$ws.Range("H107").Value = 749.1177
$ws.Range("I107").Value = 585.8333
$ws.Range("J107").Value = 1141
$ws.Range("K107").Value = 585.8333
$ws.Range("L107").Value = 1141
$ws.Range("M107").Value = 1334.1667
$ws.Range("N107").Value = -4981
# Row 122
$ws.Range("H122").Value = 4662.8
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4662.8
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 13988.4
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -18888.4
# Row 136
$ws.Range("H136").Value = 3113.492
$ws.Range("I136").Value = 1783.02
$ws.Range("K136").Value = 5349.059999999999
$ws.Range("M136").Value = -2799.059999999999

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 118.25
$ws.Range("J38").Value = 149.28572
$ws.Range("L38").Value = 447.85716
$ws.Range("N38").Value = -1141.85716
# Row 131
$ws.Range("H131").Value = 697.4400000000001
$ws.Range("J131").Value = 807.275
$ws.Range("L131").Value = 2421.825
$ws.Range("N131").Value = -12501.825

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 41668816
$ws.Range("I80").Value = 50001980
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 50001980
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -50000982
$ws.Range("N80").Value = -4996
# Row 83
$ws.Range("H83").Value = 41668816
$ws.Range("I83").Value = 50001980
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 250009900
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -250004908
$ws.Range("N83").Value = -24984
# Row 113
$ws.Range("H113").Value = 1408.3846
$ws.Range("I113").Value = 1508.9412
$ws.Range("J113").Value = 1218.4445
$ws.Range("K113").Value = 1508.9412
$ws.Range("L113").Value = 1218.4445
$ws.Range("M113").Value = 661.0588
$ws.Range("N113").Value = -5558.4445
# Row 126
$ws.Range("H126").Value = 3246.95
$ws.Range("I126").Value = 2817.705
$ws.Range("J126").Value = 4768.8184
$ws.Range("K126").Value = 8453.115
$ws.Range("L126").Value = 14306.4552
$ws.Range("M126").Value = -5983.115
$ws.Range("N126").Value = -19246.4552

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2300.6
$ws.Range("I16").Value = 2167
$ws.Range("J16").Value = 2501
$ws.Range("K16").Value = 2167
$ws.Range("L16").Value = 2501
$ws.Range("M16").Value = -1997
$ws.Range("N16").Value = -2841
# Row 136
$ws.Range("H136").Value = 4803.684
$ws.Range("I136").Value = 1692
$ws.Range("J136").Value = 7604.2
$ws.Range("K136").Value = 5076
$ws.Range("L136").Value = 22812.6
$ws.Range("M136").Value = -2526
$ws.Range("N136").Value = -27912.6

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 39900
$ws.Range("J80").Value = 39900
$ws.Range("L80").Value = 39900
$ws.Range("N80").Value = -41896
# Row 83
$ws.Range("H83").Value = 39900
$ws.Range("J83").Value = 39900
$ws.Range("L83").Value = 119700
$ws.Range("N83").Value = -129684
# Row 136
$ws.Range("H136").Value = 10197.8
$ws.Range("I136").Value = 5026.25
$ws.Range("J136").Value = 13645.5
$ws.Range("K136").Value = 15078.75
$ws.Range("L136").Value = 40936.5
$ws.Range("M136").Value = -12528.75
$ws.Range("N136").Value = -46036.5
